$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Seguimiento a riesgos")

# Fill in row 20 with the new risk entry
$ws.Range("A20").Value = 10
$ws.Range("B20").Value = "Los integrantes del equipo de desarrollo no administran de manera correcta el sofware de gestion de proyecto"
$ws.Range("C20").Value = 4
$ws.Range("D20").Value = 4
$ws.Range("E20").Formula = "=D20*C20"
$ws.Range("F20").Value = "utilizar y actualizar el software de gestion para administrar el proyecto a diario"
$ws.Range("G20").Value = "Cambiar de software para llevar la gestion del seguimiento del proyecto"
$ws.Range("H20").Value = (Get-Date -Year 2019 -Month 3 -Day 5)
$ws.Range("I20").Value = "Presentado"

$ws.Rows.Item(20).RowHeight = 29.25

# Update the active selection to H20
$ws.Range("H20").Select()
